$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 changes from "[ignore]" to the birthday timestamp string.
# Prefix with an apostrophe so Excel stores it as quote-prefixed text
# (matching the cell's existing text style) instead of reinterpreting it.
$ws.Range("D3").Value = "'2015-02-11 13:22:11"

# Update the active cell / selection shown in the sheet view.
$ws.Range("B18").Select()
